# "Generate Report for Handoff"
# Adds a new row (row 4) describing the handoff of
# f48566f5-e7cf-497f-9071-cdc761036d89.md to each of the three report
# sheets (Overview, zh-cn, de-de) and grows each sheet's table/autofilter
# range + dimension accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "f48566f5-e7cf-497f-9071-cdc761036d89.md"
$wsOverview.Range("B4").Value = "e2e\f48566f5-e7cf-497f-9071-cdc761036d89.md"
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/293157402b5922e87ff62584654fbc4c8b00c2fc/e2e/f48566f5-e7cf-497f-9071-cdc761036d89.md",
    "",
    "",
    "e2e\f48566f5-e7cf-497f-9071-cdc761036d89.md"
) | Out-Null
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-09-09 13:12:57"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A4").Value = "f48566f5-e7cf-497f-9071-cdc761036d89.md"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8c9507fabcd24dc4074237a2f690a945aa15e5bb/e2e/f48566f5-e7cf-497f-9071-cdc761036d89.md",
    "",
    "",
    "f48566f5-e7cf-497f-9071-cdc761036d89.md"
) | Out-Null
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "False"
$wsZhCn.Range("G4").Value = "f48566f5-e7cf-497f-9071-cdc761036d89.8c9507fabcd24dc4074237a2f690a945aa15e5bb.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-09-09 13:12:46"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value = ""
$wsZhCn.Range("J4").Value = ""
$wsZhCn.Range("K4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A4").Value = "f48566f5-e7cf-497f-9071-cdc761036d89.md"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8c9507fabcd24dc4074237a2f690a945aa15e5bb/e2e/f48566f5-e7cf-497f-9071-cdc761036d89.md",
    "",
    "",
    "f48566f5-e7cf-497f-9071-cdc761036d89.md"
) | Out-Null
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "False"
$wsDeDe.Range("G4").Value = "f48566f5-e7cf-497f-9071-cdc761036d89.8c9507fabcd24dc4074237a2f690a945aa15e5bb.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-09-09 13:12:57"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value = ""
$wsDeDe.Range("J4").Value = ""
$wsDeDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""
